$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 242, pushing existing rows 242:332 down to 243:333
$ws.Rows.Item(242).Insert()

# Populate the newly inserted row 242 with the new data record
$ws.Range("A242").Value = 10
$ws.Range("B242").Value = "Vega Modelo de Temuco"
$ws.Range("C242").Value = "La Araucanía"
$ws.Range("D242").Value = 44636
$ws.Range("E242").Value = 9
$ws.Range("F242").Value = 100112037
$ws.Range("G242").Value = "Cebollín"
$ws.Range("H242").Value = "Sin especificar"
$ws.Range("I242").Value = "Primera"
$ws.Range("J242").Value = 20
$ws.Range("K242").Value = 8000
$ws.Range("L242").Value = 8000
$ws.Range("M242").Value = 8000
$ws.Range("N242").Value = "$/docena de paquetes"
$ws.Range("O242").Value = "Provincia de Cautín"
$ws.Range("P242").Value = 667
$ws.Range("Q242").Value = 12
$ws.Range("R242").Value = "Hortaliza"
